$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 43 (1-based), shifting existing rows 43:183 down to 44:184
$ws.Rows.Item(43).Insert()

# Populate the newly inserted row with the new entry
$ws.Cells.Item(43, 1).Value = "FOJ78948"
$ws.Cells.Item(43, 2).Value = "z4ZP-wYr2"
